# New crime data collected — weekly CompStat update (123rd Precinct)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: volume/issue number and the reporting week date range
# ---------------------------------------------------------------------------
$cell = $ws.Range("A8")
$full = $cell.Text
$idx = $full.IndexOf("49") + 1
$cell.Characters($idx, 2).Text = "50"

$cell = $ws.Range("C9")
$full = $cell.Text
$idx = $full.IndexOf("12/5/2022") + 1
$cell.Characters($idx, 9).Text = "12/12/2022"
$full = $cell.Text
$idx = $full.IndexOf("12/11/2022") + 1
$cell.Characters($idx, 10).Text = "12/18/2022"

# ---------------------------------------------------------------------------
# Helper: convert a "***.*"/"0" placeholder (text, style 14) cell into a real
# number. Copies number-format+font from a template cell of the target style
# first so the style index is reused rather than a new one minted, then
# writes the value.
# ---------------------------------------------------------------------------
function Set-NumFromText($addrInt, $addrPct, $intVal, $pctVal) {
    if ($addrInt -ne $null) {
        $ws.Range("I14").Copy() | Out-Null
        $ws.Range($addrInt).PasteSpecial(-4122) | Out-Null
        $ws.Range($addrInt).Value = $intVal
    }
    if ($addrPct -ne $null) {
        $ws.Range("M14").Copy() | Out-Null
        $ws.Range($addrPct).PasteSpecial(-4122) | Out-Null
        $ws.Range($addrPct).Value = $pctVal
    }
}

# ---------------------------------------------------------------------------
# Helper: convert a numeric cell into a text placeholder ("0" or "***.*"),
# reusing style 14 (General/text) from a stable template cell (C14).
# Order matters: write the text value FIRST (quoted, so it lands as a shared
# string) and THEN paste the template's format over it — pasting formats
# after the value keeps the string but drops the "quote prefix" flag that
# Excel would otherwise bake into a brand-new style.
# ---------------------------------------------------------------------------
function Set-TextFromNum($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range("C14").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
Set-NumFromText "D14" "E14" 1 -100
Set-NumFromText "G14" "H14" 1 -100
Set-NumFromText "J14" "K14" 1 100

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
Set-TextFromNum "C16" "0"
$ws.Range("N16").Value = -68.75

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 1
Set-TextFromNum "D17" "0"
Set-TextFromNum "E17" "***.*"
$ws.Range("I17").Value = 50
$ws.Range("K17").Value = 61.290322580645
$ws.Range("L17").Value = -20.634920634920
$ws.Range("M17").Value = 2.040816326530
$ws.Range("N17").Value = -52.830188679245

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
Set-NumFromText "C18" $null 1 $null
Set-NumFromText "D18" "E18" 1 0
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = 400
$ws.Range("I18").Value = 42
$ws.Range("J18").Value = 24
$ws.Range("K18").Value = 75
$ws.Range("L18").Value = -2.325581395348
$ws.Range("M18").Value = -59.615384615384
$ws.Range("N18").Value = -87.155963302752

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -40
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -43.243243243243
$ws.Range("I19").Value = 262
$ws.Range("J19").Value = 195
$ws.Range("K19").Value = 34.358974358974
$ws.Range("L19").Value = 48.022598870056
$ws.Range("M19").Value = 80.689655172413
$ws.Range("N19").Value = 11.965811965812

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -11.111111111111
$ws.Range("I20").Value = 119
$ws.Range("J20").Value = 59
$ws.Range("K20").Value = 101.694915254237
$ws.Range("L20").Value = 170.454545454545
$ws.Range("M20").Value = 221.621621621622
$ws.Range("N20").Value = -82.902298850574

# ---------------------------------------------------------------------------
# Row 21 - TOTAL (bold styles; values only, style unchanged)
# ---------------------------------------------------------------------------
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = -21.428571428571
$ws.Range("F21").Value = 43
$ws.Range("G21").Value = 50
$ws.Range("H21").Value = -14
$ws.Range("I21").Value = 497
$ws.Range("J21").Value = 322
$ws.Range("K21").Value = 54.347826086956
$ws.Range("L21").Value = 44.476744186046
$ws.Range("M21").Value = 37.673130193905
$ws.Range("N21").Value = -65.293296089385

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 44
$ws.Range("G24").Value = 38
$ws.Range("H24").Value = 15.789473684210
$ws.Range("I24").Value = 481
$ws.Range("J24").Value = 275
$ws.Range("K24").Value = 74.909090909090
$ws.Range("L24").Value = 74.909090909090
$ws.Range("M24").Value = -8.901515151515

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 178
$ws.Range("J25").Value = 140
$ws.Range("K25").Value = 27.142857142857
$ws.Range("L25").Value = 35.877862595419
$ws.Range("M25").Value = -16.822429906542

# ---------------------------------------------------------------------------
# Row 28 - Shooting Vic.
# ---------------------------------------------------------------------------
Set-NumFromText "D28" "E28" 1 -100
Set-NumFromText "G28" "H28" 1 -100
$ws.Range("J28").Value = 2
$ws.Range("K28").Value = 50
Set-NumFromText $null "M28" $null 200

# ---------------------------------------------------------------------------
# Row 29 - Shooting Inc.
# ---------------------------------------------------------------------------
Set-NumFromText "D29" "E29" 1 -100
Set-NumFromText "G29" "H29" 1 -100
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = 50
Set-NumFromText $null "M29" $null 200

Write-Output "edit complete"
